$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Target values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).
# This reflects the weekly data refresh described in the commit message -
# rows keep their identity (market/category/etc.) but the date/price/volume
# figures are reshuffled among rows.
$updates = @(
    @{Row=2;  D=44498; J=40; K=4000; L=4000; M=4000; P=4000},
    @{Row=3;  D=44504; J=55; K=4000; L=4000; M=4000; P=4000},
    @{Row=4;  D=44749; J=65; K=6000; L=6000; M=6000; P=6000},
    @{Row=5;  D=44508; J=30; K=4000; L=4000; M=4000; P=4000},
    @{Row=6;  D=44781; J=40; K=5000; L=5000; M=5000; P=5000},
    @{Row=7;  D=44680; J=20; K=5000; L=5000; M=5000; P=5000},
    @{Row=8;  D=44312; J=50; K=4000; L=4000; M=4000; P=4000},
    @{Row=9;  D=44509; J=20; K=4000; L=4000; M=4000; P=4000},
    @{Row=10; D=44365; J=55; K=5000; L=5000; M=5000; P=5000},
    @{Row=11; D=44497; J=20; K=4000; L=4000; M=4000; P=4000},
    @{Row=12; D=44390; J=55; K=6000; L=6000; M=6000; P=6000},
    @{Row=13; D=44316; J=20; K=4000; L=4000; M=4000; P=4000},
    @{Row=14; D=44679; J=50; K=5000; L=5000; M=5000; P=5000},
    @{Row=15; D=44313; J=20; K=4000; L=4000; M=4000; P=4000},
    @{Row=16; D=44301; J=40; K=3000; L=3000; M=3000; P=3000},
    @{Row=17; D=44176; J=10; K=4000; L=4000; M=4000; P=4000},
    @{Row=18; D=44649; J=20; K=5000; L=5000; M=5000; P=5000},
    @{Row=19; D=44315; J=40; K=4000; L=4000; M=4000; P=4000},
    @{Row=20; D=44656; J=85; K=5000; L=5000; M=5000; P=5000},
    @{Row=21; D=44777; J=25; K=5000; L=5000; M=5000; P=5000},
    @{Row=22; D=44280; J=55; K=4000; L=4000; M=4000; P=4000},
    @{Row=23; D=44291; J=35; K=4000; L=4000; M=4000; P=4000},
    @{Row=24; D=44259; J=30; K=4000; L=4000; M=4000; P=4000}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value = $u.D    # D - Fecha
    $ws.Cells.Item($r, 10).Value = $u.J   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $u.K   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $u.L   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $u.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $u.P   # P - Precio $/Kg
}

$wb.Save()
